$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)

# --- Update existing data row on sheet "o_10" ---
$promptText = @'
 Given is the adjacency matrix for a weighted directed graph containing 16 nodes labelled A to P. The value corresponding to each row M and column N represents the cost of travelling between the two nodes, where 0 means no connection.   

Consider some examples

Example 1: what is the least cost path from node A to node P?
   A B C D E F G H I J K L M N O P
 A 0 2 0 0 3 0 0 0 0 0 0 0 0 0 0 0
 B 0 0 3 0 0 0 0 0 0 0 0 0 0 0 0 0
 C 0 0 0 2 0 0 3 0 0 0 0 0 0 0 0 0
 D 0 0 0 0 0 0 0 1 0 0 0 0 0 0 0 0
 E 0 0 0 0 0 1 0 0 4 0 0 0 0 0 0 0
 F 0 1 0 0 1 0 0 0 0 0 0 0 0 0 0 0
 G 0 0 0 0 0 0 0 3 0 0 4 0 0 0 0 0
 H 0 0 0 0 0 0 3 0 0 0 0 3 0 0 0 0
 I 0 0 0 0 0 0 0 0 0 0 0 0 5 0 0 0
 J 0 0 0 0 0 4 0 0 4 0 0 0 0 0 0 0
 K 0 0 0 0 0 0 0 0 0 3 0 0 0 0 2 0
 L 0 0 0 0 0 0 0 0 0 0 5 0 0 0 0 2
 M 0 0 0 0 0 0 0 0 0 0 0 0 0 3 0 0
 N 0 0 0 0 0 0 0 0 0 0 0 0 0 0 5 0
 O 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 4
 P 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0

Solution: A -> B -> C -> D -> H -> L -> P
        

Example 2: what is the least cost path from node A to node P?
   A B C D E F G H I J K L M N O P
 A 0 1 0 0 0 0 0 0 0 0 0 0 0 0 0 0
 B 0 0 0 0 0 3 0 0 0 0 0 0 0 0 0 0
 C 0 4 0 0 0 0 0 0 0 0 0 0 0 0 0 0
 D 0 0 4 0 0 0 0 0 0 0 0 0 0 0 0 0
 E 5 0 0 0 0 0 0 0 3 0 0 0 0 0 0 0
 F 0 0 0 0 1 0 0 0 0 1 0 0 0 0 0 0
 G 0 0 0 0 0 3 0 0 0 0 0 0 0 0 0 0
 H 0 0 0 3 0 0 3 0 0 0 0 0 0 0 0 0
 I 0 0 0 0 0 0 0 0 0 0 0 0 3 0 0 0
 J 0 0 0 0 0 0 0 0 0 0 4 0 0 0 0 0
 K 0 0 0 0 0 0 1 0 0 0 0 4 0 0 0 0
 L 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 5
 M 0 0 0 0 0 0 0 0 0 0 0 0 0 3 0 0
 N 0 0 0 0 0 0 0 0 0 4 0 0 0 0 1 0
 O 0 0 0 0 0 0 0 0 0 0 4 0 0 0 0 0
 P 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0

Solution: A -> B -> F -> J -> K -> L -> P
        

Example 3: what is the least cost path from node A to node I?
   A B C D E F G H I
 A 0 1 0 1 0 0 0 0 0
 B 0 0 5 0 0 0 0 0 0
 C 0 0 0 0 0 1 0 0 0
 D 0 0 0 0 0 0 1 0 0
 E 0 5 0 2 0 0 0 0 0
 F 0 0 0 0 1 0 0 0 1
 G 0 0 0 0 0 0 0 2 0
 H 0 0 0 0 4 0 0 0 5
 I 0 0 0 0 0 0 0 0 0

Solution: A -> B -> C -> F -> I
        
 Given these examples, answer the following quesiton.

what is the least cost path from node A to node P?

   A B C D E F G H I J K L M N O P
 A 0 2 0 0 4 0 0 0 0 0 0 0 0 0 0 0
 B 0 0 5 0 0 0 0 0 0 0 0 0 0 0 0 0
 C 0 0 0 0 0 0 1 0 0 0 0 0 0 0 0 0
 D 0 0 4 0 0 0 0 0 0 0 0 0 0 0 0 0
 E 0 0 0 0 0 0 0 0 3 0 0 0 0 0 0 0
 F 0 5 0 0 1 0 4 0 0 0 0 0 0 0 0 0
 G 0 0 0 0 0 0 0 2 0 0 0 0 0 0 0 0
 H 0 0 0 4 0 0 0 0 0 0 0 3 0 0 0 0
 I 0 0 0 0 0 0 0 0 0 1 0 0 1 0 0 0
 J 0 0 0 0 0 5 0 0 0 0 4 0 0 0 0 0
 K 0 0 0 0 0 0 4 0 0 0 0 0 0 0 4 0
 L 0 0 0 0 0 0 0 0 0 0 2 0 0 0 0 0
 M 0 0 0 0 0 0 0 0 0 0 0 0 0 5 0 0
 N 0 0 0 0 0 0 0 0 0 3 0 0 0 0 5 0
 O 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 4
 P 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0
    
'@

$ws1.Range("A2").Value = $promptText
$ws1.Range("B2").Value = "A -> E -> I -> J -> K -> O -> P"
$ws1.Range("C2").Value = "The least cost path from node A to node P is A -> E -> J -> K -> O -> P."
# D2 ("Wrong") is unchanged

# --- Add new column E: evaluator_partial_correctness ---
$ws1.Range("D1").Copy($ws1.Range("E1"))
$ws1.Range("E1").Value = "evaluator_partial_correctness"
$ws1.Range("E2").Value = "Output: 5/7"

# --- Add new sheets o_20 and o_20_jumbled, each with the same header row ---
$ws2 = $wb.Worksheets.Add($null, $ws1)
$ws2.Name = "o_20"
$ws3 = $wb.Worksheets.Add($null, $ws2)
$ws3.Name = "o_20_jumbled"

foreach ($ws in @($ws2, $ws3)) {
    $ws1.Range("A1:E1").Copy($ws.Range("A1"))
}

$wb.Worksheets.Item(1).Activate()

